# Update the "Förändrad" (changed) date column (C) for data rows 2-14
# from serial date 45208 (2023-10-09) to 45212 (2023-10-13).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 3).Value = 45212
}
